$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 'Dom Cobb and Arthur are "extractors" who performcorporate espionageusing experimental dream-sharing technology to infiltrate their targets''subconsciousand extract information. Their latest target, Saito, is impressed with Cobb''s ability tolayer multiple dreams within each other. He offers to hire Cobb for the ostensibly impossible job of implanting an idea into a person''s subconscious; performing "inception" on Robert Fischer, the son of Saito''s competitor Maurice Fischer, with the idea to dissolve his father''s company. In return, Saito promises to clear Cobb''s criminal status, allowing him to return home to his children.'
$ws.Range("D2").Value = 'Christopher Nolan'
$ws.Range("E2").Value = 'Leonardo DiCaprio, Leonardo DiCaprio, Ken Watanabe, Ken Watanabe, Joseph Gordon-Levitt, Joseph Gordon-Levitt, Marion Cotillard, Marion Cotillard, Elliot Page[a], Elliot Page, [a], Tom Hardy, Tom Hardy, Cillian Murphy, Cillian Murphy, Tom Berenger, Tom Berenger, Michael Caine, Michael Caine'
$ws.Range("F2").Value = 'July 8, 2010(2010-07-08)(Odeon Leicester Square)July 16, 2010(2010-07-16)(United States and United Kingdom)'
$ws.Range("G2").Value = '148 minutes[1]'
$ws.Range("C3").Value = 'In 1945, theNew York CityCorleone familydon,Vito Corleone, listens to requests during his daughterConnie''s wedding toCarlo Rizzi. Vito''s youngest sonMichael, aMarineandWorld War IIhero who has so far stayed out of the family business, introduces his girlfriendKay Adamsto his family at the reception. Johnny Fontane, a popular singer and Vito''sgodson, seeks Vito''s help in securing a movie role. Vito sends hisconsigliereTom Hagento persuade studio president Jack Woltz to offer Johnny the part. Woltz refuses Hagen''s request at first, but soon complies after finding the severed head of his prized stud horse in his bed.'
$ws.Range("D3").Value = 'Francis Ford Coppola'
$ws.Range("E3").Value = 'Marlon Brando, Marlon Brando, Al Pacino, Al Pacino, James Caan, James Caan, Richard Castellano, Richard Castellano, Robert Duvall, Robert Duvall, Sterling Hayden, Sterling Hayden, John Marley, John Marley, Richard Conte, Richard Conte, Diane Keaton, Diane Keaton'
$ws.Range("F3").Value = 'March 14, 1972(1972-03-14)(Loew''s State Theatre)March 24, 1972(1972-03-24)(United States)'
$ws.Range("G3").Value = '175 minutes[1]'
$ws.Range("C4").Value = 'In the mid-21st century, humanity facesextinctiondue todust stormsand widespread cropblights. Joseph Cooper, a widowed formerNASAtest pilot, works as a farmer and raises his children, Murph and Tom, alongside his father-in-law Donald. Living in apost-truth society, Cooper is reprimanded by Murph''s teachers for telling her that theApollo missionswere not fabricated. During a dust storm, the two discover that dust patterns in Murph''s room, which she first attributes to aghost, result from a gravitational anomaly, and translate intogeographic coordinates. These lead them to a secret NASA facility headed by Professor John Brand, who explains that, 48 years earlier, awormholeappeared nearSaturn, leading to a system in anothergalaxywith twelve potentially habitable planets located near ablack holenamed Gargantua. Volunteers of theLazarusexpedition had previously travelled through the wormhole to evaluate the planets, with Miller, Edmunds, and Mann reporting back desirable results.'
$ws.Range("D4").Value = 'Christopher Nolan'
$ws.Range("E4").Value = 'Matthew McConaughey, Matthew McConaughey, Anne Hathaway, Anne Hathaway, Jessica Chastain, Jessica Chastain, Bill Irwin, Bill Irwin, Ellen Burstyn, Ellen Burstyn, Michael Caine, Michael Caine'
$ws.Range("F4").Value = 'October 26, 2014(2014-10-26)(TCL Chinese Theatre)November 5, 2014(2014-11-05)(United States)November 7, 2014(2014-11-07)(United Kingdom)'
$ws.Range("G4").Value = '169 minutes[1]'
$ws.Range("C5").Value = 'In 1981, a feather lands at a bus stop inSavannah, Georgia, where a man named Forrest Gump recounts his life story to strangers who sit next to him on a bench.'
$ws.Range("D5").Value = 'Robert Zemeckis'
$ws.Range("E5").Value = 'Tom Hanks, Tom Hanks, Robin Wright, Robin Wright, Gary Sinise, Gary Sinise, Mykelti Williamson, Mykelti Williamson, Sally Field, Sally Field'
$ws.Range("F5").Value = 'June 23, 1994(1994-06-23)(Los Angeles)July 6, 1994(1994-07-06)(United States)'
$ws.Range("G5").Value = '142 minutes'
$ws.Range("C6").Value = 'In 2018, 23 days afterThanoserased half of all life in the universe,[a]Carol DanversrescuesTony StarkandNebulafrom deep space. They reunite with the remainingAvengers—Bruce Banner,Steve Rogers,Thor,Natasha Romanoff, andJames Rhodes—andRocketon Earth. Locating Thanos onan uninhabited planet, they plan to use theInfinity Stonesto reverse his actions but find that Thanos has destroyed them. Enraged, Thor decapitates Thanos.'
$ws.Range("D6").Value = 'Anthony RussoJoe Russo, Anthony Russo, Joe Russo'
$ws.Range("E6").Value = 'Robert Downey Jr., Robert Downey Jr., Chris Evans, Chris Evans, Mark Ruffalo, Mark Ruffalo, Chris Hemsworth, Chris Hemsworth, Scarlett Johansson, Scarlett Johansson, Jeremy Renner, Jeremy Renner, Don Cheadle, Don Cheadle, Paul Rudd, Paul Rudd, Brie Larson, Brie Larson, Karen Gillan, Karen Gillan, Danai Gurira, Danai Gurira, Benedict Wong, Benedict Wong, Jon Favreau, Jon Favreau, Bradley Cooper, Bradley Cooper, Gwyneth Paltrow, Gwyneth Paltrow, Josh Brolin, Josh Brolin'
$ws.Range("F6").Value = 'April 22, 2019(2019-04-22)(Los Angeles Convention Center)April 26, 2019(2019-04-26)(United States)'
$ws.Range("G6").Value = '181 minutes[1]'
